# Config.xlsx — "Assets" sheet: add six new generic-asset config rows
# (CurrencyCode_New, ExchangeRate_New, SellingRate_New, BuyingRate_New,
# UseDrawer_New, UseDrawer_EditCurrency), each mirroring the existing
# Name/Value/Asset layout, and move the sheet selection down to the new
# first empty row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assets")
$null = $ws.Select()

# Copy the formatting of the last six populated rows (22:27) down into the
# six blank rows that follow (28:33). Inserting shifts every row below
# down by six, so afterwards we drop the six blank rows that got pushed
# past the sheet's original last row to keep the sheet the same size.
$null = $ws.Range("A22:C27").Copy()
$null = $ws.Range("A28:C33").EntireRow.Insert()
$ws.Range("A28:C33").RowHeight = 14.25
$null = $ws.Range("A1000:C1005").EntireRow.Delete()

# Fill in the new asset setting names (columns A and B mirror each other,
# column C is always "Generic Asset", same as the rows above).
$newNames = @(
    "CurrencyCode_New",
    "ExchangeRate_New",
    "SellingRate_New",
    "BuyingRate_New",
    "UseDrawer_New",
    "UseDrawer_EditCurrency"
)

$row = 28
foreach ($name in $newNames) {
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $name
    $ws.Cells.Item($row, 3).Value = "Generic Asset"
    $row = $row + 1
}

# Move the selection to the next empty row below the new entries.
$null = $ws.Range("A34").Select()
